$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '58.075.75'
$ws.Range("E2").Value = '  -4.42%  '
# Row 3
$ws.Range("D3").Value = '2.618.43'
$ws.Range("E3").Value = '  -2.98%  '
# Row 4
$ws.Range("E4").Value = '  +0.12%  '
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '516.28'
$ws.Range("E5").Value = '  -1.79%  '
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '142.29'
$ws.Range("E6").Value = '  -1.99%  '
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.25%  '
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.567'
$ws.Range("E8").Value = '  -1.76%  '
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '6.72'
$ws.Range("E9").Value = '  +0.37%  '
# Row 10
$ws.Range("E10").Value = '  -3.08%  '
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.336'
$ws.Range("E11").Value = '  -1.05%  '
# Row 12
$ws.Range("E12").Value = '  +1.26%  '
# Row 13
$ws.Range("D13").Value = '3.078.73'
$ws.Range("E13").Value = '  -3.01%  '
# Row 14
$ws.Range("D14").Value = '58.090.08'
$ws.Range("E14").Value = '  -4.04%  '
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '20.64'
$ws.Range("E15").Value = '  -2.91%  '
# Row 16
$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000135'
$ws.Range("E16").Value = '  -1.82%  '
# Row 17
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '2.612.67'
$ws.Range("E17").Value = '  -8.14%  '
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '4.40'
$ws.Range("E18").Value = '  -2.44%  '
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '334.75'
$ws.Range("E19").Value = '  -3.15%  '
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.33'
$ws.Range("E20").Value = '  -2.81%  '
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.25'
$ws.Range("E21").Value = '  -3.07%  '
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.00'
$ws.Range("E22").Value = '  +0.16%  '
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '63.75'
$ws.Range("E23").Value = '  +0.35%  '
# Row 24
$ws.Range("E24").Value = '  +0.01%  '
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.166'
$ws.Range("E25").Value = '  -2.55%  '
# Row 26
$ws.Range("E26").Value = '  +0.71%  '
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.05'
$ws.Range("E27").Value = '  -3.27%  '
# Row 28
$ws.Range("D28").Value = '0.0₃0781'
$ws.Range("E28").Value = '  -4.57%  '
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.59'
$ws.Range("E29").Value = '  -3.14%  '
# Row 30
$ws.Range("E30").Value = '  +0.11%  '
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.57'
$ws.Range("E31").Value = '  -1.46%  '
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '151.14'
$ws.Range("E32").Value = '  +0.79%  '
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '18.68'
$ws.Range("E33").Value = '  -1.91%  '
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.07'
$ws.Range("E34").Value = '  -4.10%  '
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.17'
$ws.Range("E35").Value = '  -5.31%  '
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.890'
$ws.Range("E36").Value = '  -5.59%  '
# Row 37
$ws.Range("E37").Value = '  -1.44%  '
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.841'
$ws.Range("E38").Value = '  -3.50%  '
# Row 39
$ws.Range("E39").Value = '  -5.98%  '
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.60'
$ws.Range("E40").Value = '  -1.79%  '
# Row 41
$ws.Range("E41").Value = '  +0.36%  '
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.598'
$ws.Range("E42").Value = '  -1.83%  '
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0964'
$ws.Range("E43").Value = '  -2.27%  '
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '267.79'
$ws.Range("E44").Value = '  -5.33%  '
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.61'
$ws.Range("E45").Value = '  +1.30%  '
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '19.11'
$ws.Range("E46").Value = '  -4.95%  '
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0531'
$ws.Range("E47").Value = '  -1.36%  '
# Row 48
$ws.Range("D48").Value = '2.030.65'
$ws.Range("E48").Value = '  -5.23%  '
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0227'
$ws.Range("E49").Value = '  -2.33%  '
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '4.62'
$ws.Range("E50").Value = '  -3.55%  '
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '18.20'
$ws.Range("E51").Value = '  -4.75%  '

Write-Output "Applied cryptos list update."
